# FEAT #2134 ajouter balise pour code barres chrono
# Add a new NAME/VAR row pointing to the chrono bar-code attachment tag,
# along with the help text to insert the barcode image above the tag.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New row 48: tag name + its descriptive label
$ws.Range("A48").Value = "[attachments.chronoBarCode;ope=changepic]"
$ws.Range("B48").Value = "Insérer une image juste au dessus de cette balise"

# Column A needs to be a bit wider to comfortably show the new (longer) tag
$ws.Columns.Item(1).ColumnWidth = 42.66

# Scroll the view down so the newly added row is visible, and select it
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A48:B48").Select()
